$d = $word.ActiveDocument

$d.Content.Find.Execute("64×23=", $true, $false, $false, $false, $false, $true, 1, $false, "64×47=", 2) | Out-Null
$d.Content.Find.Execute("39×61=", $true, $false, $false, $false, $false, $true, 1, $false, "71×84=", 2) | Out-Null
$d.Content.Find.Execute("87×68=", $true, $false, $false, $false, $false, $true, 1, $false, "83×58=", 2) | Out-Null
$d.Content.Find.Execute("90×36=", $true, $false, $false, $false, $false, $true, 1, $false, "65×13=", 2) | Out-Null
$d.Content.Find.Execute("29×23=", $true, $false, $false, $false, $false, $true, 1, $false, "94×18=", 2) | Out-Null
$d.Content.Find.Execute("91×74=", $true, $false, $false, $false, $false, $true, 1, $false, "29×16=", 2) | Out-Null
$d.Content.Find.Execute("45×42=", $true, $false, $false, $false, $false, $true, 1, $false, "86×24=", 2) | Out-Null
$d.Content.Find.Execute("86×32=", $true, $false, $false, $false, $false, $true, 1, $false, "95×78=", 2) | Out-Null
$d.Content.Find.Execute("36×48=", $true, $false, $false, $false, $false, $true, 1, $false, "18×79=", 2) | Out-Null
$d.Content.Find.Execute("60×57=", $true, $false, $false, $false, $false, $true, 1, $false, "22×80=", 2) | Out-Null
$d.Content.Find.Execute("40×98=", $true, $false, $false, $false, $false, $true, 1, $false, "69×23=", 2) | Out-Null
$d.Content.Find.Execute("67×91=", $true, $false, $false, $false, $false, $true, 1, $false, "57×34=", 2) | Out-Null
$d.Content.Find.Execute("54×28=", $true, $false, $false, $false, $false, $true, 1, $false, "55×59=", 2) | Out-Null
$d.Content.Find.Execute("81×21=", $true, $false, $false, $false, $false, $true, 1, $false, "38×67=", 2) | Out-Null
$d.Content.Find.Execute("19×99=", $true, $false, $false, $false, $false, $true, 1, $false, "45×83=", 2) | Out-Null
$d.Content.Find.Execute("85×18=", $true, $false, $false, $false, $false, $true, 1, $false, "51×21=", 2) | Out-Null
$d.Content.Find.Execute("97×61=", $true, $false, $false, $false, $false, $true, 1, $false, "14×81=", 2) | Out-Null
$d.Content.Find.Execute("52×86=", $true, $false, $false, $false, $false, $true, 1, $false, "87×77=", 2) | Out-Null
$d.Content.Find.Execute("91×86=", $true, $false, $false, $false, $false, $true, 1, $false, "93×76=", 2) | Out-Null
$d.Content.Find.Execute("57×78=", $true, $false, $false, $false, $false, $true, 1, $false, "89×25=", 2) | Out-Null
$d.Content.Find.Execute("79×99=", $true, $false, $false, $false, $false, $true, 1, $false, "47×74=", 2) | Out-Null
$d.Content.Find.Execute("76×59=", $true, $false, $false, $false, $false, $true, 1, $false, "96×95=", 2) | Out-Null
$d.Content.Find.Execute("26×50=", $true, $false, $false, $false, $false, $true, 1, $false, "40×82=", 2) | Out-Null
$d.Content.Find.Execute("65×48=", $true, $false, $false, $false, $false, $true, 1, $false, "72×91=", 2) | Out-Null
$d.Content.Find.Execute("80×84=", $true, $false, $false, $false, $false, $true, 1, $false, "59×93=", 2) | Out-Null
